$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The three model blocks (gemma_7b, llama3.2_3b, mistral_7b), each spanning 8 rows
# (zero_shot/Raw, Role-based, one_shot/Raw, Role-based, few_shot/Raw, Role-based,
# auto_cot/Raw, Role-based) get cyclically rotated among rows 18-25, 26-33, 34-41.
# New order (top to bottom): mistral_7b, gemma_7b, llama3.2_3b.
# Only column A labels (at the first row of each block) and the D:G numeric
# values change; columns B and C keep their existing pattern.

# Update the model name labels (merged cells A18:A25, A26:A33, A34:A41)
$ws.Range("A18").Value = "mistral_7b"
$ws.Range("A26").Value = "gemma_7b"
$ws.Range("A34").Value = "llama3.2_3b"

# New values for D18:G41 (precision, recall, f1-score, accuracy)
$data = @(
    @(0.83, 0.2, 0.32, 0.76),
    @(0.8100000000000001, 0.28, 0.42, 0.78),
    @(0.72, 0.43, 0.54, 0.79),
    @(0.7, 0.49, 0.58, 0.8),
    @(0.71, 0.45, 0.55, 0.79),
    @(0.7, 0.48, 0.57, 0.8),
    @(0.8100000000000001, 0.24, 0.37, 0.77),
    @(0.84, 0.23, 0.36, 0.77),
    @(0.44, 0.59, 0.5, 0.67),
    @(0.67, 0.5, 0.57, 0.79),
    @(0.39, 0.74, 0.51, 0.6),
    @(0.44, 0.8, 0.57, 0.66),
    @(0.4, 0.71, 0.52, 0.62),
    @(0.53, 0.61, 0.57, 0.74),
    @(0.43, 0.67, 0.53, 0.66),
    @(0.62, 0.52, 0.57, 0.78),
    @(0.43, 0.6899999999999999, 0.53, 0.66),
    @(0.54, 0.67, 0.6, 0.75),
    @(0.5, 0.66, 0.57, 0.72),
    @(0.46, 0.8, 0.59, 0.68),
    @(0.66, 0.36, 0.46, 0.77),
    @(0.6, 0.62, 0.61, 0.78),
    @(0.4, 0.77, 0.53, 0.62),
    @(0.52, 0.7, 0.59, 0.73)
)

$startRow = 18
for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Range("D$r").Value = $row[0]
    $ws.Range("E$r").Value = $row[1]
    $ws.Range("F$r").Value = $row[2]
    $ws.Range("G$r").Value = $row[3]
}

$wb.Save()
